# Updates the "想去人数" (interest-count) figures in column F that were
# refreshed by the gh-pages data-publishing job (commit 456a3b4).
# Sheets touched: 展览 (Exhibitions), 演出 (Performances), 全部类型 (All types).
# 本地生活 (Local life) is untouched by this refresh.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F5").Value = 1002
$wsExhibition.Range("F6").Value = 5592
$wsExhibition.Range("F7").Value = 504
$wsExhibition.Range("F8").Value = 707
$wsExhibition.Range("F9").Value = 965
$wsExhibition.Range("F13").Value = 592
$wsExhibition.Range("F15").Value = 23
$wsExhibition.Range("F17").Value = 1881
$wsExhibition.Range("F18").Value = 1481
$wsExhibition.Range("F19").Value = 942
$wsExhibition.Range("F21").Value = 199
$wsExhibition.Range("F23").Value = 562
$wsExhibition.Range("F24").Value = 161
$wsExhibition.Range("F25").Value = 1057
$wsExhibition.Range("F28").Value = 3016
$wsExhibition.Range("F30").Value = 106
$wsExhibition.Range("F32").Value = 129
$wsExhibition.Range("F33").Value = 39
$wsExhibition.Range("F34").Value = 411
$wsExhibition.Range("F37").Value = 15
$wsExhibition.Range("F40").Value = 742
$wsExhibition.Range("F42").Value = 55

$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F4").Value = 205
$wsPerformance.Range("F6").Value = 143

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F5").Value = 1002
$wsAllTypes.Range("F7").Value = 5592
$wsAllTypes.Range("F8").Value = 504
$wsAllTypes.Range("F9").Value = 707
$wsAllTypes.Range("F11").Value = 205
$wsAllTypes.Range("F12").Value = 965
$wsAllTypes.Range("F15").Value = 143
$wsAllTypes.Range("F18").Value = 592
$wsAllTypes.Range("F20").Value = 23
$wsAllTypes.Range("F23").Value = 1881
$wsAllTypes.Range("F24").Value = 1481
$wsAllTypes.Range("F25").Value = 942
$wsAllTypes.Range("F26").Value = 199
$wsAllTypes.Range("F29").Value = 562
$wsAllTypes.Range("F30").Value = 161
$wsAllTypes.Range("F31").Value = 1057
$wsAllTypes.Range("F32").Value = 3016
$wsAllTypes.Range("F34").Value = 106
$wsAllTypes.Range("F36").Value = 129
$wsAllTypes.Range("F37").Value = 39
$wsAllTypes.Range("F38").Value = 411
$wsAllTypes.Range("F41").Value = 15
$wsAllTypes.Range("F43").Value = 742

